# cryptos.xlsx price/volume refresh (GitHub Actions scheduled update).
#
# Columns D (Price) and E (Volume(1h)) are stored as plain text in the
# source sheet (note the dotted-thousands formatting, e.g. "57.845.50",
# and the padded "  +1.23%  " volume strings), so every write below goes
# through Range.Value as a string.
#
# Caveat: Excel auto-detects plain decimals (e.g. "1.00", "375.92") typed
# into a General-formatted cell and silently coerces them to numbers,
# which would drop the trailing zero / exact text. For any such D-column
# value we therefore set NumberFormat to Text ("@") immediately before
# assigning it, exactly as a spreadsheet author would do by hand. Values
# that can never be mis-parsed as numbers (two-dot thousands separators,
# the "₃3"-subscript ticker, percent strings with padding/spaces, plain
# text) are left with their existing General format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.908.08'
$ws.Range('E2').Value = '  -1.00%  '
$ws.Range('D3').Value = '2.936.36'
$ws.Range('E3').Value = '  +0.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '550.28'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.76'
$ws.Range('E6').Value = '  +9.38%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.31%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.513'
$ws.Range('E8').Value = '  +4.13%  '
$ws.Range('D9').Value = '2.930.32'
$ws.Range('E9').Value = '  +1.00%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.127'
$ws.Range('E10').Value = '  +0.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '4.77'
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('E12').Value = '  +3.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000221'
$ws.Range('E13').Value = '  +3.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.92'
$ws.Range('E14').Value = '  +4.46%  '
$ws.Range('E15').Value = '  +2.52%  '
$ws.Range('D16').Value = '3.414.99'
$ws.Range('E16').Value = '  +0.56%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.88'
$ws.Range('E17').Value = '  +7.85%  '
$ws.Range('D18').Value = '2.928.00'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('D19').Value = '57.831.48'
$ws.Range('E19').Value = '  -1.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '416.76'
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.32'
$ws.Range('E21').Value = '  +3.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.695'
$ws.Range('E22').Value = '  +6.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.41'
$ws.Range('E23').Value = '  +6.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.99'
$ws.Range('E24').Value = '  +2.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '79.48'
$ws.Range('E25').Value = '  +2.86%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('E28').Value = '  +0.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.03'
$ws.Range('E29').Value = '  +6.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.45'
$ws.Range('E30').Value = '  +6.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.27'
$ws.Range('E31').Value = '  +2.92%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.99'
$ws.Range('E32').Value = '  +0.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0977'
$ws.Range('E33').Value = '  +3.78%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.68'
$ws.Range('E34').Value = '  +4.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.939'
$ws.Range('E35').Value = '  +4.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.08'
$ws.Range('E36').Value = '  +6.23%  '
$ws.Range('D37').Value = '0.0₃0697'
$ws.Range('E37').Value = '  +12.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.31'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.73'
$ws.Range('E39').Value = '  +4.93%  '
$ws.Range('E40').Value = '  +10.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.108'
$ws.Range('E41').Value = '  +2.23%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '375.92'
$ws.Range('E42').Value = '  +6.29%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0345'
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').Value = '2.692.42'
$ws.Range('E44').Value = '  +2.93%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '123.64'
$ws.Range('E46').Value = '  +4.46%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.237'
$ws.Range('E47').Value = '  +4.40%  '
$ws.Range('E48').Value = '  +2.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.97'
$ws.Range('E49').Value = '  +2.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.99'
$ws.Range('E50').Value = '  +1.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.00'
$ws.Range('E51').Value = '  +2.46%  '
